$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel
# (pure decimal-looking price values) are forced to remain text, matching the
# original inline-string cell type, then the style is restored to the default
# "Normal" so no stray formatting is introduced.

$ws.Range('D2').Value = '97.755.85'
$ws.Range('E2').Value = '  -1.20%  '
$ws.Range('D3').Value = '3.414.08'
$ws.Range('E3').Value = '  +3.11%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '256.52'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '658.97'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.47%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.46'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.431'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +5.57%  '
$ws.Range('E9').Value = '  +9.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.999'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Value = '3.412.74'
$ws.Range('E11').Value = '  +3.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.214'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +6.83%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '42.21'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.90%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.50'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +18.72%  '
$ws.Range('E15').Value = '  +3.48%  '
$ws.Range('D16').Value = '97.397.61'
$ws.Range('E16').Value = '  -1.26%  '
$ws.Range('D17').Value = '4.048.20'
$ws.Range('E17').Value = '  +3.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.56'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +35.73%  '
$ws.Range('D19').Value = '3.405.03'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.67'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +12.75%  '
$ws.Range('E21').Value = '  +71.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.89'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +14.80%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.46'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '509.40'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000205'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.60%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.16'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +9.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '96.28'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +8.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.73'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.45%  '
$ws.Range('D29').Value = '3.563.61'
$ws.Range('E29').Value = '  +2.21%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.152'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +12.32%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '11.35'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +12.65%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.995'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.52%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.195'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.37%  '
$ws.Range('E34').Value = '  +0.47%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.573'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +21.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '29.77'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.51%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.19'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +12.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.82'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +8.49%  '
$ws.Range('E39').Value = '  +5.88%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '515.12'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.13%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.39'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +12.14%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '24.71'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.56%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.853'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +7.95%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0431'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +29.91%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.66'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.17%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.33'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.73%  '
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.42'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +14.67%  '
$ws.Range('B48').Value = 'Cosmos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.23'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +12.25%  '
$ws.Range('B49').Value = 'USDe'
$ws.Range('C49').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.00'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.10'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.65%  '
$ws.Range('B51').Value = 'ImmutableX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.57'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +15.38%  '
